$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 20191011
$ws.Range("J3").Value = 20191062
$ws.Range("J4").Value = "did not close"
$ws.Range("J5").Value = "did not close"
$ws.Range("J6").Value = "did not close"
$ws.Range("J7").Value = "did not close"
$ws.Range("J8").Value = 20190096
$ws.Range("J9").Value = "did not close"
$ws.Range("J10").Value = 20182625
$ws.Range("J11").Value = "did not close"
